# Auto-generated Excel COM-interop script
# Updates crypto price/volume data per commit "Updated cryptos list on Fri Jul 28 17:39:09 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.259.75"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.870.62"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7114"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.49"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3106"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07696"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.04"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08351"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "1.876.81"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.217"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7100"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.18"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "29.289.95"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008263"
$ws.Range("E17").Value = "  +5.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.931"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.92"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").Value = "2.131.58"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.839"
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1622"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.14"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.993"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.49"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.408"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.338"
$ws.Range("E31").Value = "  +5.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.277"
$ws.Range("E32").Value = "  -5.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05244"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.922"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7491"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.171"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.680"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01857"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.714"
$ws.Range("D40").Value = "1.152.17"
$ws.Range("E40").Value = "  -1.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.358"
$ws.Range("E41").Value = "  +4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.00"
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8869"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.53"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "2.027.98"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5189"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.792"
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000120"
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.358"
$ws.Range("E50").Value = "  +0.61%  "

Write-Output "Updated cryptos list"
